$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 283 (existing rows 283:321 shift down
# to 285:323). Excel copies formatting from the row above, which keeps the
# date-formatted style on column D.
$ws.Rows("283:284").Insert()

# --- New row 283 ---
$ws.Cells.Item(283, 1).Value = 4
$ws.Cells.Item(283, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(283, 3).Value = "Los Lagos"
$ws.Cells.Item(283, 4).Value = 44491
$ws.Cells.Item(283, 5).Value = 10
$ws.Cells.Item(283, 6).Value = 100112004
$ws.Cells.Item(283, 7).Value = "Cebolla"
$ws.Cells.Item(283, 8).Value = "Morada(o)"
$ws.Cells.Item(283, 9).Value = "1a (cosecha)"
$ws.Cells.Item(283, 10).Value = 200
$ws.Cells.Item(283, 11).Value = 11000
$ws.Cells.Item(283, 12).Value = 11000
$ws.Cells.Item(283, 13).Value = 11000
$ws.Cells.Item(283, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(283, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(283, 16).Value = 611
$ws.Cells.Item(283, 17).Value = 18
$ws.Cells.Item(283, 18).Value = "Hortaliza"

# --- New row 284 ---
$ws.Cells.Item(284, 1).Value = 4
$ws.Cells.Item(284, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(284, 3).Value = "Los Lagos"
$ws.Cells.Item(284, 4).Value = 44491
$ws.Cells.Item(284, 5).Value = 10
$ws.Cells.Item(284, 6).Value = 100112004
$ws.Cells.Item(284, 7).Value = "Cebolla"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "1a (cosecha)"
$ws.Cells.Item(284, 10).Value = 900
$ws.Cells.Item(284, 11).Value = 7500
$ws.Cells.Item(284, 12).Value = 8000
$ws.Cells.Item(284, 13).Value = 7750
$ws.Cells.Item(284, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(284, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(284, 16).Value = 431
$ws.Cells.Item(284, 17).Value = 18
$ws.Cells.Item(284, 18).Value = "Hortaliza"
